$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "without_project"
$ws.Range("B4").Value = 240129

$ws.Rows.Item(3).RowHeight = 13.8
$ws.Rows.Item(4).RowHeight = 13.8

$ws.Range("B4").Select()
